$d = $word.ActiveDocument

# "Tallow for glass" -> "Sandever"
$d.Content.Find.Execute("Tallow for glass", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sandever", 2)

# "lightens well" -> "lighten well" (second occurrence of "lightens" in the sentence)
$d.Content.Find.Execute("All things that lightens well", $true, $false, $false, $false, $false,
                         $true, 1, $false, "All things that lighten well", 2)
